$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/Volume columns (D/E) hold numeric-looking text values; force Text format
# per-cell so Excel keeps the exact literal string instead of parsing it as a number/percent.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.59"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.29%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "29.64"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.16%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.318"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.46%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05742"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.67%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.653"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.98%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.204"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.92%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8583"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.03%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8522"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.95%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1381"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.91%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07105"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.06%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03211"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "11.84%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09328"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.64%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001547"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.38%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005949"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.90%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.514"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.88%"
$ws.Range("B17").Value = "BTSEToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.196"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.05%"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.01013"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1,594.61%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3161"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.51%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03356"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.65%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.29%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.493"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "20.09%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04139"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.42%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1408"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2.04%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001223"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.30%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004167"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-18.11%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001200"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.81%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001450"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-25.22%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.77%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1069"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.18%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.002417"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "15.10%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002951"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-48.51%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009192"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.40%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005271"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.01%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.03%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.08093"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "13.97%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002201"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-17.55%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.03%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.03%"
